# day2 - API and write into excel
# Update row 5 (Nadia Ningtias/Sidoarjo -> Nadia/surabaya) and remove the
# trailing rows (6-8: Rahmad ksmrdn, Novita, Linda) that were appended
# from the API, shrinking the sheet back down to A1:C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the data in row 5
$ws.Range("B5").Value = "Nadia"
$ws.Range("C5").Value = "surabaya"

# Drop the extra appended rows 6-8 entirely (shifts sheetData/dimension up)
$ws.Range("A6:C8").Delete()
